$d = $word.ActiveDocument
$d.Content.Find.Execute("A number of", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Several", 2)
